# Natmi following Dr Hou advice
# Rebuild the LR-pairs table: add the "ECs" sending-cluster block (rows 2-5)
# and refresh every recomputed numeric column for the full 4x4 sender/target matrix.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) is unchanged.

# Each inner array is one data row: A..T (Sending cluster .. Edge total expression derived specificity)
$data = @(
  @("ECs","Tslp","Crlf2","ECs",2,0.6666666666666666,0.5078013333333333,1.523404,0.1375609480357072,0.1375609480357073,3,1,5.179207666666667,15.537623,0.2482357269660917,0.2482357269660917,2.630008558743556,23.670077028692,0.03414754193778854,0.03414754193778855),
  @("ECs","Tslp","Crlf2","FAPs",2,0.6666666666666666,0.5078013333333333,1.523404,0.1375609480357072,0.1375609480357073,3,1,4.184307,12.552921,0.2005508480919455,0.2005508480919455,2.124796673676,19.123170063084,0.02758796479289314,0.02758796479289314),
  @("ECs","Tslp","Crlf2","M2",2,0.6666666666666666,0.5078013333333333,1.523404,0.1375609480357072,0.1375609480357073,3,1,10.70673666666667,32.12021,0.513166246835409,0.513166246835409,5.436895154982222,48.93205639484,0.07059163541460461,0.07059163541460463),
  @("ECs","Tslp","Crlf2","sCs",2,0.6666666666666666,0.5078013333333333,1.523404,0.1375609480357072,0.1375609480357073,3,1,0.793819,2.381457,0.03804717810655387,0.03804717810655387,0.4031023466253333,3.627921119628,0.005233805890420955,0.005233805890420956),
  @("FAPs","Tslp","Crlf2","ECs",3,1,1.159282,3.477846,0.3140439390222109,0.314043939022211,3,1,5.179207666666667,15.537623,0.2482357269660917,0.2482357269660917,6.004162222228667,54.037460000058,0.07795692550247349,0.07795692550247352),
  @("FAPs","Tslp","Crlf2","FAPs",3,1,1.159282,3.477846,0.3140439390222109,0.314043939022211,3,1,4.184307,12.552921,0.2005508480919455,0.2005508480919455,4.850791787574,43.657126088166,0.06298177830903963,0.06298177830903964),
  @("FAPs","Tslp","Crlf2","M2",3,1,1.159282,3.477846,0.3140439390222109,0.314043939022211,3,1,10.70673666666667,32.12021,0.513166246835409,0.513166246835409,12.41212709640667,111.70914386766,0.161156749529436,0.1611567495294361),
  @("FAPs","Tslp","Crlf2","sCs",3,1,1.159282,3.477846,0.3140439390222109,0.314043939022211,3,1,0.793819,2.381457,0.03804717810655387,0.03804717810655387,0.920260077958,8.282340701622001,0.0119484856812618,0.01194848568126181),
  @("M2","Tslp","Crlf2","ECs",3,1,1.601072333333333,4.803217,0.433722823454071,0.4337228234540711,3,1,5.179207666666667,15.537623,0.2482357269660917,0.2482357269660917,8.29228610368789,74.630574933191,0.1076655003819072,0.1076655003819072),
  @("M2","Tslp","Crlf2","FAPs",3,1,1.601072333333333,4.803217,0.433722823454071,0.4337228234540711,3,1,4.184307,12.552921,0.2005508480919455,0.2005508480919455,6.699378171873001,60.294403546857,0.08698348008054711,0.08698348008054711),
  @("M2","Tslp","Crlf2","M2",3,1,1.601072333333333,4.803217,0.433722823454071,0.4337228234540711,3,1,10.70673666666667,32.12021,0.513166246835409,0.513166246835409,17.14225985728556,154.28033871557,0.2225719134787824,0.2225719134787824),
  @("M2","Tslp","Crlf2","sCs",3,1,1.601072333333333,4.803217,0.433722823454071,0.4337228234540711,3,1,0.793819,2.381457,0.03804717810655387,0.03804717810655387,1.270961638574333,11.438654747169,0.01650192951283446,0.01650192951283446),
  @("sCs","Tslp","Crlf2","ECs",3,1,0.4233086666666666,1.269926,0.1146722894880107,0.1146722894880108,3,1,5.179207666666667,15.537623,0.2482357269660917,0.2482357269660917,2.192403491766444,19.731631425898,0.02846575914392246,0.02846575914392246),
  @("sCs","Tslp","Crlf2","FAPs",3,1,0.4233086666666666,1.269926,0.1146722894880107,0.1146722894880108,3,1,4.184307,12.552921,0.2005508480919455,0.2005508480919455,1.771253417094,15.941280753846,0.02299762490946565,0.02299762490946565),
  @("sCs","Tslp","Crlf2","M2",3,1,0.4233086666666666,1.269926,0.1146722894880107,0.1146722894880108,3,1,10.70673666666667,32.12021,0.513166246835409,0.513166246835409,4.532254422717777,40.79028980446,0.058845948412586,0.05884594841258601),
  @("sCs","Tslp","Crlf2","sCs",3,1,0.4233086666666666,1.269926,0.1146722894880107,0.1146722894880108,3,1,0.793819,2.381457,0.03804717810655387,0.03804717810655387,0.3360304624646667,3.024274162182,0.00436295702203665,0.004362957022036651)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $i + 2   # data starts at row 2
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}
